# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45989

$ws.Range("B2").Value  = 91.7
$ws.Range("C2").Value  = 87.08
$ws.Range("D2").Value  = 90.70999999999999
$ws.Range("E2").Value  = 88.89
$ws.Range("F2").Value  = 88.8
$ws.Range("G2").Value  = 92.19
$ws.Range("H2").Value  = 97.48
$ws.Range("I2").Value  = 128.66
$ws.Range("J2").Value  = 134.42
$ws.Range("K2").Value  = 101.97
$ws.Range("L2").Value  = 88.04000000000001
$ws.Range("M2").Value  = 77.89
$ws.Range("N2").Value  = 77.34999999999999
$ws.Range("O2").Value  = 71.75
$ws.Range("P2").Value  = 69.05
$ws.Range("Q2").Value  = 74.73
$ws.Range("R2").Value  = 87.04000000000001
$ws.Range("S2").Value  = 105.39
$ws.Range("T2").Value  = 120.02
$ws.Range("U2").Value  = 128.19
$ws.Range("V2").Value  = 143.07
$ws.Range("W2").Value  = 123.49
$ws.Range("X2").Value  = 102.59
$ws.Range("Y2").Value  = 99.04000000000001
$ws.Range("Z2").Value  = 98.73

$ws.Range("AB2").Value = 117.05
$ws.Range("AD2").Value = 133.28
$ws.Range("AF2").Value = 124.1
$ws.Range("AG2").Value = "0h-16h"
